# mega-import.xlsx hotfix: add "General price" / "Pricing method" / "Standard
# price" columns to the #saledata.product sheet (sheet17 / rId17).
#
# Target layout (1-indexed columns):
#   ... Q = Thue ban (unchanged)
#   R (NEW) = Gia chung / -General price cua san pham / 100000,200000,300000
#   S (was R) = Don vi tinh ton kho (unchanged content, shifted right)
#   T (NEW) = Phuong thuc dinh gia / long note / 0,0,0
#   U (NEW) = Gia tieu chuan / (blank) / 100000,200000,300000
#   V (was S) = Don vi tinh mua hang (unchanged content, shifted right)
#   W (was T) = Thue mua (unchanged content, shifted right)
#   X (was U) = Nguon cung cho viec mua (unchanged content, shifted right)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#saledata.product")
$ws.Activate()

# --- Insert the new columns -------------------------------------------------
# Step 1: one new column right before the old "R" (Don vi tinh ton kho).
$ws.Columns("R").Insert()
# Step 2: two new columns right before the old "S" (now sitting at "T" after
# the previous insert), i.e. before "Don vi tinh mua hang".
$ws.Columns("T:U").Insert()

# --- Headers / descriptions --------------------------------------------------
# Written in the same left-to-right, top-to-bottom order as the source
# workbook so new shared-string entries land at the same indices (450..454).
$ws.Cells.Item(1, 18).Value = "Giá chung"
$ws.Cells.Item(2, 18).Value = "-General price của sản phẩm"
$ws.Cells.Item(1, 20).Value = "Phương thức định giá"
$ws.Cells.Item(1, 21).Value = "Giá tiêu chuẩn"
$ws.Cells.Item(2, 20).Value = "-0 là nhập trước xuất trước`n-1 là bình quân gia quyền`n-2 là thực tế đích danh`n- Bắt buộc nếu sản phẩm lưu kho"
# Cell U2 (col 21) stays blank, matching the source diff.

# --- Data rows 3-5 (R = Gia chung, T = Phuong thuc dinh gia, U = Gia tieu chuan) ---
$prices = @(100000, 200000, 300000)
for ($i = 0; $i -lt 3; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 18).Value = $prices[$i]
    $ws.Cells.Item($row, 20).Value = 0
    $ws.Cells.Item($row, 21).Value = $prices[$i]
}

# Column insertion copies the format of the column to the left, which already
# gives R/S/U the correct style. Column T inherits the "no-wrap" style from
# its left neighbour (S) instead of the wrapping style used by the rest of
# the data rows, so fix that up explicitly to match the target formatting.
$ws.Range("T3:T5").WrapText = $true
$ws.Range("T3:T5").VerticalAlignment = -4160

# --- View bookkeeping (best effort; matches the author's selection) ---------
$ws.Range("T9").Select()

$wb.Save()
